$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bibliography paragraph style: add alignment "left" (core of the commit)
# ---------------------------------------------------------------------------
$biblioStyle = $d.Styles("Bibliography")
$biblioStyle.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft

# ---------------------------------------------------------------------------
# 2) Replace the "List of references." placeholder paragraph with the two
#    real bibliography entries.
# ---------------------------------------------------------------------------
$refsPara = $d.Paragraphs.Last

# --- second bibliography entry (reuse the existing paragraph) -------------
$text2 = "Vellutini BC and Hejnol A (2016). Expression of segment polarity genes in brachiopods supports a non-segmental ancestral role of engrailed for bilaterians. https://doi.org/10.1038/srep32387"
$refsPara.Range.Text = $text2

$refsPara.Format.Style = "Bibliography"
$refsPara.Format.SpaceBefore = 9
$refsPara.Format.SpaceAfter = 9

$p2start = $refsPara.Range.Start

# whole-paragraph defaults common to every run in entry 2
$whole2 = $d.Range($p2start, $p2start + 189)
$whole2.Font.Size = 12
$whole2.Font.Color = 3617564   # RGB(0x37,0x39,0x3C) -> BGR int for wdColor
$whole2.Font.Spacing = 0

# run 1: "Vellutini BC and Hejnol A (2016). " (plain)
$r2a = $d.Range($p2start + 0, $p2start + 34)
$r2a.Font.Bold = 0
$r2a.Font.Italic = 0

# run 2: "Expression of ... role of " (bold)
$r2b = $d.Range($p2start + 34, $p2start + 129)
$r2b.Font.Bold = 1
$r2b.Font.Italic = 0

# run 3: "engrailed" (bold italic)
$r2c = $d.Range($p2start + 129, $p2start + 138)
$r2c.Font.Bold = 1
$r2c.Font.Italic = 1

# run 4: " for bilaterians" (bold)
$r2d = $d.Range($p2start + 138, $p2start + 154)
$r2d.Font.Bold = 1
$r2d.Font.Italic = 0

# run 5: ". " (plain)
$r2e = $d.Range($p2start + 154, $p2start + 156)
$r2e.Font.Bold = 0
$r2e.Font.Italic = 0

# run 6: the DOI link, turned into a real hyperlink with the InternetLink style
$r2f = $d.Range($p2start + 156, $p2start + 189)
$r2f.Font.Bold = 0
$r2f.Font.Italic = 0
$d.Hyperlinks.Add($r2f, "https://doi.org/10.1038/srep32387") | Out-Null
$r2fAfter = $d.Range($p2start + 156, $p2start + 189)
$r2fAfter.Style = "InternetLink"

# --- first bibliography entry (brand-new paragraph inserted before) -------
$refsPara = $d.Paragraphs.Last
$refsPara.Range.InsertParagraphBefore()
$entry1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$text1 = "Vellutini BC, Martín-Durán JM, and Hejnol A (2017). Cleavage modification did not alter blastomere fates during bryozoan evolution. https://doi.org/10.1186/s12915-017-0371-9"
$entry1.Range.Text = $text1
$entry1.Format.Style = "Bibliography"

$p1start = $entry1.Range.Start

$whole1 = $d.Range($p1start, $p1start + 173)
$whole1.Font.Size = 12

# run 1: "Vellutini BC, Martín-Durán JM, and Hejnol A (2017). " (plain)
$r1a = $d.Range($p1start + 0, $p1start + 52)
$r1a.Font.Bold = 0

# run 2: title (bold)
$r1b = $d.Range($p1start + 52, $p1start + 130)
$r1b.Font.Bold = 1

# run 3: ". " (plain)
$r1c = $d.Range($p1start + 130, $p1start + 132)
$r1c.Font.Bold = 0

# run 4: DOI link
$r1d = $d.Range($p1start + 132, $p1start + 173)
$r1d.Font.Bold = 0
$d.Hyperlinks.Add($r1d, "https://doi.org/10.1186/s12915-017-0371-9") | Out-Null
$r1dAfter = $d.Range($p1start + 132, $p1start + 173)
$r1dAfter.Style = "InternetLink"

Write-Output "bibliography updated"
